$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(8)
$shp.Left = -3760 / 12700
$shp.Top = 2299521 / 12700
